$d = $word.ActiveDocument
$r = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$xml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t>Descripción de ideas descartadas</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Idea nº1: </w:t>
      </w:r>
      <w:r>
        <w:t>U</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">ltrasonidos: </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:ind w:firstLine="696"/>
      </w:pPr>
      <w:r>
        <w:t>P</w:t>
      </w:r>
      <w:r>
        <w:t>ros: fácil implementación, bajo coste, alta fiabilidad.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:ind w:firstLine="696"/>
      </w:pPr>
      <w:r>
        <w:t>Contras: necesarias 2 puertas (entrada y salida), no hay forma de identificar a la persona.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Idea nº2: </w:t>
      </w:r>
      <w:r>
        <w:t>R</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">econocimiento facial: </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:ind w:firstLine="696"/>
      </w:pPr>
      <w:r>
        <w:t>P</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">ros: muy fiable. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:ind w:firstLine="696"/>
      </w:pPr>
      <w:r>
        <w:t>Contras: no necesarias 2 puertas (entrada y salida), alto coste de producción.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Idea nº3: </w:t>
      </w:r>
      <w:r>
        <w:t>I</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">nfrarrojos: </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:ind w:firstLine="696"/>
      </w:pPr>
      <w:r>
        <w:t>P</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">ros: fácil implementación y bajo coste. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:ind w:firstLine="696"/>
      </w:pPr>
      <w:r>
        <w:t>Contras: baja fiabilidad, 2 puertas y no hay forma de identificar.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Idea nº4: modulo NFC: </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:ind w:firstLine="696"/>
      </w:pPr>
      <w:r>
        <w:t>P</w:t>
      </w:r>
      <w:r>
        <w:t>ros: fácil implementación, bajo coste, 1 sola puerta, se puede identificar a la gente, alta fiabilidad</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
      <w:r>
        <w:t>Conclusi</w:t>
      </w:r>
      <w:r>
        <w:t>ó</w:t>
      </w:r>
      <w:r>
        <w:t>n: al principio pensamos que lo mas eficiente ser</w:t>
      </w:r>
      <w:r>
        <w:t>ía</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>un ultrasonidos</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>, pero según empezamos a desarrollar ideas nos decantamos por el NFC por todas las ventajas respecto al resto de opciones.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t>Tecnologías a utilizar</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Respecto a la tecnología a utilizar, podemos dividirla en 4 partes: </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>percepction</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>layer</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">(tecnología referida al hardware a utilizar), </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>network</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>layer</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>(tecnología referida a conexiones</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">inalámbricas), </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>processing</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>layer</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>(tecnología con la que captar la</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">información hardware) y por último </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>application</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>layer</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>(tecnología que utilizará el usuario final).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">En la primera capa, disponemos de una placa </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>NodeMCU</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> para transmitir los datos al servidor y un </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>NFID</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> que captará nuestros llaveros NFC. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
      <w:r>
        <w:t>Para la segunda capa usaremos el wifi</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">de 2,4 </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Ghz</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> que lleva la propia placa</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
      <w:r>
        <w:t>P</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">ara la capa de proceso utilizaremos un servidor </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>M</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>qtt</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> para recibir la información de la placa, un servidor apache y una base de datos</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>M</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>ysql</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
      </w:pPr>
      <w:r>
        <w:t>Por último, en la última capa</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">haremos uso de </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>R</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>eact</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> para crear un entorno web y posiblemente utilizaremos una aplicación Android.</w:t>
      </w:r>
    </w:p>
'@
$r.InsertXML($xml)
